$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.338.59"
$ws.Range("E2").Value = "  +0.81%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.149.12"
$ws.Range("E3").Value = "  +1.18%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.69"
$ws.Range("E5").Value = "  -1.44%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "619.93"
$ws.Range("E6").Value = "  -1.26%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.13"
$ws.Range("E7").Value = "  +1.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.381"
$ws.Range("E8").Value = "  +2.50%  "

# Row 9
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.147.44"
$ws.Range("E10").Value = "  +14.31%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.748"
$ws.Range("E11").Value = "  +0.98%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.204"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  -0.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.23"
$ws.Range("E14").Value = "  -1.29%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.60"
$ws.Range("E15").Value = "  +1.56%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.106.74"
$ws.Range("E16").Value = "  +0.74%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.733.86"
$ws.Range("E17").Value = "  +1.54%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.154.40"
$ws.Range("E18").Value = "  +1.91%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.73"
$ws.Range("E19").Value = "  -3.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.05"
$ws.Range("E20").Value = "  +4.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.02"
$ws.Range("E21").Value = "  +4.21%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "458.82"
$ws.Range("E22").Value = "  +2.39%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000203"
$ws.Range("E23").Value = "  -3.81%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.18"
$ws.Range("E24").Value = "  +1.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.95"
$ws.Range("E25").Value = "  -1.11%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.63"
$ws.Range("E26").Value = "  +61.77%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "89.13"
$ws.Range("E27").Value = "  -4.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.86"
$ws.Range("E28").Value = "  -2.70%  "

# Row 29
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.153"
$ws.Range("E29").Value = "  +33.91%  "

# Row 30
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.318.32"
$ws.Range("E30").Value = "  +1.48%  "

# Row 31
$ws.Range("E31").Value = "  -0.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.236"
$ws.Range("E32").Value = "  +9.72%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.41"
$ws.Range("E34").Value = "  +0.63%  "

# Row 35
$ws.Range("E35").Value = "  +13.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.47"
$ws.Range("E36").Value = "  -0.30%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.51"
$ws.Range("E37").Value = "  -1.97%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("E38").Value = "  +1.74%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "493.23"
$ws.Range("E39").Value = "  -0.68%  "

# Row 40
$ws.Range("B40").Value = "MantraDAO"
$ws.Range("C40").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.90"
$ws.Range("E40").Value = "  -12.11%  "

# Row 41
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.33"
$ws.Range("E41").Value = "  +1.74%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.448"
$ws.Range("E42").Value = "  +7.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.36"
$ws.Range("E43").Value = "  -8.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.13"
$ws.Range("E44").Value = "  +0.04%  "

# Row 45
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.715"
$ws.Range("E46").Value = "  +3.55%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.93"
$ws.Range("E47").Value = "  +0.23%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "155.71"
$ws.Range("E48").Value = "  -2.66%  "

# Row 49
$ws.Range("E49").Value = "  +1.15%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.52"
$ws.Range("E50").Value = "  -1.02%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0328"
$ws.Range("E51").Value = "  +5.80%  "
